$d = $word.ActiveDocument

$d.Content.Find.Execute("Chemistry in St. Andrews", $true, $false, $false, $false, $false, $true, 1, $false, "Philosophy in Glasgow", 2)
$d.Content.Find.Execute("in chemistry", $true, $false, $false, $false, $false, $true, 1, $false, "in philosophy", 2)
$d.Content.Find.Execute("worked in Agfa", $true, $false, $false, $false, $false, $true, 1, $false, "worked in KODAK", 2)
